$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to text first so numeric-looking values like
# "1.001" or "8.320" are stored as literal strings (matching the source data)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.429.51"
$ws.Range("E2").Value = "  +3.58%  "

$ws.Range("D3").Value = "1.868.13"
$ws.Range("E3").Value = "  +1.98%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "337.69"
$ws.Range("E5").Value = "  +1.99%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "0.4685"
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("D8").Value = "0.3972"
$ws.Range("E8").Value = "  +3.60%  "

$ws.Range("D9").Value = "47.75"
$ws.Range("E9").Value = "  +2.53%  "

$ws.Range("E10").Value = "  +1.76%  "

$ws.Range("D11").Value = "0.9985"
$ws.Range("E11").Value = "  +2.72%  "

$ws.Range("D12").Value = "21.96"
$ws.Range("E12").Value = "  +4.15%  "

$ws.Range("D13").Value = "6.046"
$ws.Range("E13").Value = "  +2.74%  "

$ws.Range("D14").Value = "1.865.48"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").Value = "7.262"
$ws.Range("E15").Value = "  +2.95%  "

$ws.Range("D16").Value = "90.65"
$ws.Range("E16").Value = "  +3.03%  "

$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("D19").Value = "0.06619"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "17.53"
$ws.Range("E20").Value = "  +2.19%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "28.443.09"
$ws.Range("E22").Value = "  +3.59%  "

$ws.Range("E23").Value = "  +2.53%  "

$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("E25").Value = "  -1.40%  "

$ws.Range("D26").Value = "2.084.64"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").Value = "160.76"
$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").Value = "19.77"

$ws.Range("D29").Value = "2.117"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("D30").Value = "5.487"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").Value = "120.04"
$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("D32").Value = "0.9702"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("D33").Value = "0.09511"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").Value = "3.588"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("D36").Value = "1.374"
$ws.Range("E36").Value = "  +4.54%  "

$ws.Range("D37").Value = "0.06094"
$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("D38").Value = "0.02248"
$ws.Range("E38").Value = "  +2.24%  "

$ws.Range("D39").Value = "8.320"
$ws.Range("E39").Value = "  +3.41%  "

$ws.Range("D40").Value = "1.181"
$ws.Range("E40").Value = "  +2.51%  "

$ws.Range("D41").Value = "0.5932"
$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "0.1875"
$ws.Range("E43").Value = "  +1.91%  "

$ws.Range("D44").Value = "10.30"

$ws.Range("D45").Value = "1.281"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D46").Value = "0.5557"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").Value = "12.16"
$ws.Range("E47").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +4.47%  "

$ws.Range("D49").Value = "0.07252"
$ws.Range("E49").Value = "  +9.15%  "

$ws.Range("D50").Value = "2.065"
$ws.Range("E50").Value = "  +13.52%  "

$ws.Range("D51").Value = "111.82"
$ws.Range("E51").Value = "  +1.28%  "

# Restore the default (unstyled) cell style so the Price column keeps its
# original look (no lingering custom text-format style index).
$ws.Range("D2:D51").Style = "Normal"
